# GW24 league table update
# Week 23 (row 21) results for Eren (B), Mert (C), Arda (D)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = 99
$ws.Range("C21").Value = 106
$ws.Range("D21").Value = 72

# Match the light "carried-over" formatting used by the row above (row 20)
$ws.Range("B20:D20").Copy()
$ws.Range("B21:D21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move/restore the active cell selection to N22, as recorded after the edit
$ws.Range("N22").Select()
